# Natmi following Dr Hou advice
# Update Vegfa-Nrp2 LR-pair stats: ligand/receptor-expressing cell counts
# changed from 1 to 3, which propagates into total-expression and
# specificity columns (G:J, M:P) and edge-weight columns (Q:T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ E=3; G=36.35426266666666;  H=109.062788;        I=0.4094848412143908; J=0.4094848412143908; K=3; M=26.83081766666667; N=80.49245300000001; O=0.5916656861001716; P=0.5916656861001716; Q=975.4145930154406;  R=8778.731337138965;  S=0.2422781295247323;  T=0.2422781295247323  }
  3  = @{ E=3; G=36.35426266666666;  H=109.062788;        I=0.4094848412143908; J=0.4094848412143908; K=3; M=10.21969166666667; N=30.659075;         O=0.2253617819930474; P=0.2253617819930474; Q=371.5293552223445;  R=3343.7641970011;    S=0.09228223351521518; T=0.09228223351521518 }
  4  = @{ E=3; G=36.35426266666666;  H=109.062788;        I=0.4094848412143908; J=0.4094848412143908; K=3; M=8.297426666666667; N=24.89228;          O=0.1829725319067811; P=0.1829725319067811; Q=301.6468284974044;  R=2714.82145647664;   S=0.0749244781744433;  T=0.0749244781744433  }
  5  = @{ E=3; G=45.11545066666667;  H=135.346352;        I=0.5081685556916724; J=0.5081685556916724; K=3; M=26.83081766666667; N=80.49245300000001; O=0.5916656861001716; P=0.5916656861001716; Q=1210.484430786829;  R=10894.35987708146;  S=0.3006658971578466;  T=0.3006658971578466  }
  6  = @{ E=3; G=45.11545066666667;  H=135.346352;        I=0.5081685556916724; J=0.5081685556916724; K=3; M=10.21969166666667; N=30.659075;         O=0.2253617819930474; P=0.2253617819930474; Q=461.0659952160446;  R=4149.593956944401;  S=0.1145217712635084;  T=0.1145217712635084  }
  7  = @{ E=3; G=45.11545066666667;  H=135.346352;        I=0.5081685556916724; J=0.5081685556916724; K=3; M=8.297426666666667; N=24.89228;          O=0.1829725319067811; P=0.1829725319067811; Q=374.3421434402845;  R=3369.079290962561;  S=0.09298088727031738; T=0.09298088727031738 }
  8  = @{ E=3; G=7.310771333333334;  H=21.932314;         I=0.0823466030939367; J=0.0823466030939367; K=3; M=26.83081766666667; N=80.49245300000001; O=0.5916656861001716; P=0.5916656861001716; Q=196.1539726473603;  R=1765.385753826242;  S=0.04872165941759257; T=0.04872165941759257 }
  9  = @{ E=3; G=7.310771333333334;  H=21.932314;         I=0.0823466030939367; J=0.0823466030939367; K=3; M=10.21969166666667; N=30.659075;         O=0.2253617819930474; P=0.2253617819930474; Q=74.71382887217224;  R=672.4244598495501;  S=0.01855777721432376; T=0.01855777721432376 }
  10 = @{ E=3; G=7.310771333333334;  H=21.932314;         I=0.0823466030939367; J=0.0823466030939367; K=3; M=8.297426666666667; N=24.89228;          O=0.1829725319067811; P=0.1829725319067811; Q=60.66058901510223;  R=545.94530113592;    S=0.01506716646202037; T=0.01506716646202037 }
}

foreach ($rowNum in $data.Keys) {
  $cols = $data[$rowNum]
  foreach ($colLetter in $cols.Keys) {
    $addr = "$colLetter$rowNum"
    $ws.Range($addr).Value = $cols[$colLetter]
  }
}
